# Apply the workbook changes described by the diff:
#  1. Insert a new "Player Info" sheet as the first sheet with player
#     biographical data (ID, NAME, BATTING_HAND, BOWL_STYLE).
#  2. In the existing "ODI Batting" sheet, rename the MATCH_CARD_LINK
#     column to MATCH_CODE and replace the full scorecard URLs with just
#     the numeric match code.
#  3. In the existing "ODI Bowling" sheet, make the same MATCH_CARD_LINK
#     -> MATCH_CODE header rename and URL -> numeric code replacement.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "Player Info" sheet and move it to the front so the
#    resulting sheet order is: Player Info, ODI Batting, ODI Bowling.
# ---------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"
$playerInfo.Move($wb.Worksheets.Item(1))

# Header row (bold, thin border, centered horizontally, top vertically -
# matching the header style already used on the other sheets).
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data row - force text storage (to match the other sheets, where every
# value - including numeric-looking ones - is stored as text) and then
# strip the resulting number-format style back off so the cell keeps the
# default "Normal" style, same as the other plain data cells.
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4201"
$playerInfo.Range("A2").Style = "Normal"

$playerInfo.Range("B2").Value = "Michael Gertges Neser"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium"

# ---------------------------------------------------------------------
# 2. "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")

$odiBatting.Range("D1").Value = "MATCH_CODE"

$odiBatting.Range("D2").NumberFormat = "@"
$odiBatting.Range("D2").Value = "4166"
$odiBatting.Range("D2").Style = "Normal"

$odiBatting.Range("D3").NumberFormat = "@"
$odiBatting.Range("D3").Value = "4169"
$odiBatting.Range("D3").Style = "Normal"

# ---------------------------------------------------------------------
# 3. "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------
$odiBowling = $wb.Worksheets.Item("ODI Bowling")

$odiBowling.Range("B1").Value = "MATCH_CODE"

$odiBowling.Range("B2").NumberFormat = "@"
$odiBowling.Range("B2").Value = "4166"
$odiBowling.Range("B2").Style = "Normal"

$odiBowling.Range("B3").NumberFormat = "@"
$odiBowling.Range("B3").Value = "4169"
$odiBowling.Range("B3").Style = "Normal"
